# "Vue et Ctrl SupprimerUtilisateur"
#
# On the "Liste des tâches" sheet (Feuil1), mark the two tasks that were
# still "À faire" as finished:
#   - row 62: "Dévelopement et test CtrlSupprimerUtilisateur"
#   - row 63: "Dévelopement et test VueSupprimerUtilisateur"
#
# For each row: clear the "A faire" (col B) mark, set "Terminé" (col D) and
# "Validé" (col H) to X, assign "Mathieu" as Responsable (col E), and stamp
# the "Date de fin" (col G) with 08/11/2016 - matching the styling already
# used by the neighbouring finished rows (e.g. row 59).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

foreach ($r in 62, 63) {
    $ws.Cells.Item($r, 2).Value = $null        # B - A faire (cleared)
    $ws.Cells.Item($r, 4).Value = "X"          # D - Terminé
    $ws.Cells.Item($r, 5).Value = "Mathieu"    # E - Responsable
    $ws.Cells.Item($r, 7).Value = 42682        # G - Date de fin (08/11/2016)
    $ws.Cells.Item($r, 8).Value = "X"          # H - Validé
}

# Reuse the date-cell number format/style already applied to the other
# "Date de fin" entries so no new style gets minted.
$ws.Range("G59").Copy()
$ws.Range("G62:G63").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the scroll position / selection captured with the edit.
$ws.Activate()
$ws.Range("H49").Select()
